$wb = $excel.ActiveWorkbook

# Reference sheet to source the existing city-row cell formatting (font Verdana 12, FF333333)
$srcSheet = $wb.Worksheets.Item("San Luis Obispo")

# Data for the new "Kern County" worksheet: Zip, City, County, AreaCode
$kernData = @(
    @("93528", "Johannesburg", "Kern", "760"),
    @("93516", "Boron", "Kern", "760"),
    @("93222", "Pine Mountain Club ", "Kern", "661"),
    @("93313", "Bakersfield", "Kern", "661"),
    @("93383", "Bakersfield", "Kern", "661"),
    @("93301", "Bakersfield", "Kern", "661"),
    @("93596", "Boron", "Kern", "760"),
    @("93384", "Bakersfield", "Kern", "661"),
    @("93287", "Woody", "Kern", "661"),
    @("93308", "Bakersfield", "Kern", "661"),
    @("93309", "Bakersfield", "Kern", "661"),
    @("93556", "Ridgecrest", "Kern", "760"),
    @("93385", "Bakersfield", "Kern", "661"),
    @("93554", "Randsburg", "Kern", "760"),
    @("93220", "Edison", "Kern", "661"),
    @("93555", "Ridgecrest", "Kern", "760"),
    @("93302", "Bakersfield", "Kern", "661"),
    @("93312", "Bakersfield", "Kern", "661"),
    @("93311", "Bakersfield", "Kern", "661"),
    @("93558", "Red Mountain", "Kern", "760"),
    @("93518", "Caliente", "Kern", "661"),
    @("93307", "Bakersfield", "Kern", "661"),
    @("93561", "Tehachapi", "Kern", "661"),
    @("93519", "Cantil", "Kern", "760"),
    @("93560", "Rosamond", "Kern", "661"),
    @("93304", "Bakersfield", "Kern", "661"),
    @("93380", "Bakersfield", "Kern", "661"),
    @("93303", "Bakersfield", "Kern", "661"),
    @("93581", "Tehachapi", "Kern", "661"),
    @("93531", "Keene", "Kern", "661"),
    @("93314", "Bakersfield", "Kern", "661"),
    @("93305", "Bakersfield", "Kern", "661"),
    @("93306", "Bakersfield", "Kern", "661"),
    @("93224", "Fellows", "Kern", "661"),
    @("93523", "Edwards", "Kern", "661"),
    @("93251", "Mc Kittrick", "Kern", "661"),
    @("93225", "Frazier Park", "Kern", "661"),
    @("93206", "Buttonwillow", "Kern", "661"),
    @("93215", "Delano", "Kern", "661"),
    @("93252", "Maricopa", "Kern", "661"),
    @("93255", "Onyx", "Kern", "760"),
    @("93527", "Inyokern", "Kern", "760"),
    @("93389", "Bakersfield", "Kern", "661"),
    @("93222", "Frazier Park", "Kern", "661"),
    @("93240", "Lake Isabella", "Kern", "760"),
    @("93238", "Kernville", "Kern", "760"),
    @("93390", "Bakersfield", "Kern", "661"),
    @("93241", "Lamont", "Kern", "661"),
    @("93203", "Arvin", "Kern", "661"),
    @("93524", "Edwards", "Kern", "661"),
    @("93250", "Mc Farland", "Kern", "661"),
    @("93205", "Bodfish", "Kern", "760"),
    @("93249", "Lost Hills", "Kern", "661"),
    @("93243", "Lebec", "Kern", "661"),
    @("93226", "Glennville", "Kern", "661"),
    @("93280", "Wasco", "Kern", "661"),
    @("93386", "Bakersfield", "Kern", "661"),
    @("93276", "Tupman", "Kern", "661"),
    @("93388", "Bakersfield", "Kern", "661"),
    @("93387", "Bakersfield", "Kern", "661"),
    @("93505", "California City", "Kern", "760"),
    @("93285", "Wofford Heights", "Kern", "760"),
    @("93283", "Weldon", "Kern", "760"),
    @("93502", "Mojave", "Kern", "661"),
    @("93504", "California City", "Kern", "760"),
    @("93501", "Mojave", "Kern", "661"),
    @("93216", "Delano", "Kern", "661"),
    @("93268", "Taft", "Kern", "661"),
    @("93263", "Shafter", "Kern", "661")

)

# Add the new sheet after the last existing sheet (San Luis Obispo)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Kern County"

$rowCount = $kernData.Count

# Copy the formatting (font/style) of the template data row onto the full A1:D<n> block
$srcSheet.Range("A1:D1").Copy($ws.Range("A1:D$rowCount"))

# Match the original row height (16) used throughout the county sheets
$ws.Range("A1:D$rowCount").RowHeight = 16

# Fill in the data. Column C is written before column B on each row so that new shared
# strings are minted in the same order as the source workbook ("Kern" before the city names).
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 1
    $row = $kernData[$i]
    $ws.Cells.Item($r, 1).Value = [int]$row[0]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = [int]$row[3]
}

# Column B needs to be wide enough to show the longest city name
$ws.Columns.Item(2).ColumnWidth = 38.166666666666664

# Leave the same selection Excel would have shown after typing the last few rows
$ws.Range("G30").Select()

Write-Output "done"
